$wb = $excel.ActiveWorkbook

# --- Rename the existing sheet "Sayfa1" -> "binek" ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "binek"

# --- Add the new "LCV" sheet right after "binek" ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "LCV"

# --- Bring over the formatting (font/number-format) from "binek" so the
#     style table is reused instead of duplicated ---
$ws1.Range("A1:B1").Copy() | Out-Null
$ws2.Range("A1:B1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws1.Range("B4:B5").Copy() | Out-Null
$ws2.Range("B4:B5").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws1.Range("B6").Copy() | Out-Null
$ws2.Range("B6").PasteSpecial(-4122) | Out-Null      # xlPasteFormats

# --- Fill in the LCV data (OTV/hurda/kredi values for LCV) ---
$ws2.Range("A1").Value = "degisken"
$ws2.Range("B1").Value = "deger"

$ws2.Range("A2").Value = "lcv_kredi_orani"
$ws2.Range("B2").Value = 0.5

$ws2.Range("A3").Value = "lcv_ortalama_vade (yil)"
$ws2.Range("B3").Value = 3

$ws2.Range("A4").Value = "mevcut_yillik faiz"
$ws2.Range("B4").Value = 0.12

$ws2.Range("A5").Value = "indirimli_yillik_faiz"
$ws2.Range("B5").Value = 0.09

$ws2.Range("A6").Value = "lcv_max indirimli kredi miktari"
$ws2.Range("B6").Value = 100000

# --- Column widths on LCV (closest achievable match to the authored sizes) ---
$ws2.Columns("A").ColumnWidth = 25.5
$ws2.Columns("B").ColumnWidth = 10.666666666666666

# --- Selections: binek ends up with A1:B6 selected (no longer the active
#     tab), LCV ends up active with B6 selected/active-cell ---
$ws1.Range("A1:B6").Select() | Out-Null
$ws2.Select() | Out-Null
$ws2.Range("B6").Select() | Out-Null
